$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two rows that were dropped from the dataset ("RM 232" and "SC 92").
# Deleting row 26 first ("RM 232") shifts "SC 92" up from row 28 to row 27,
# so it is then removed by deleting row 27.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# Toggle which rows have an imputed value present in column F.
$ws.Range("F3").Value = 17.64
$ws.Range("F5").ClearContents()
$ws.Range("F21").Value = 16.58
$ws.Range("F23").ClearContents()

# The row for "SC 193" (now row 32 after the deletions above) gains its
# previously-missing column F value.
$ws.Range("F32").Value = 17.39
